$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C7 value (was "S1AGL065", now the branch address string)
$ws.Range("C7").Value = "NUESTRA SRA.DEL BUEN VIAJE 739"

# Update the active selection shown in the saved sheet view (was C10, now D11)
$ws.Range("D11").Select()
